$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.261557579040527
$ws.Range("B1").Value = 2.325201034545898
$ws.Range("C1").Value = 3.028860807418823
$ws.Range("D1").Value = 3.4892258644104
$ws.Range("E1").Value = 1.4439697265625
